# Insert two new data-center / colo rows into Sheet1:
#   - "BOD" (Bordeaux, France) before the existing "QWJ" (Americana, Brazil) row (row 166)
#   - "SAT" (San Antonio, United States) before the existing "ADL" (Adelaide, Australia) row
#     (which, after the first insertion, sits at row 291)
# Both insertions push the rows below them down by one, matching the new
# sheet dimension A1:G303. The rest of the data (everything before row 166,
# and everything between the two new rows) is left untouched and simply
# shifts down with the native row insert.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Insert "Bordeaux, France" (BOD) at row 166 -----------------------------
$ws.Rows.Item(166).Insert()

# Match the "colo" column look used by every other data row: bold, thin box
# border, centered/top-aligned.
$ws.Range("A166").Font.Bold = $true
$ws.Range("A166").HorizontalAlignment = -4108
$ws.Range("A166").VerticalAlignment = -4160
$ws.Range("A166").Borders.LineStyle = 1

$ws.Range("A166").Value = "BOD"
$ws.Range("B166").Value = "Bordeaux, France"
$ws.Range("C166").Value = 44.82946
$ws.Range("D166").Value = -0.58355
$ws.Range("E166").Value = "FR"
$ws.Range("F166").Value = "Europe"
$ws.Range("G166").Value = "Bordeaux"

# --- Insert "San Antonio, United States" (SAT) at row 291 ------------------
# After the row-166 insertion, the former row 290 ("ADL" / Adelaide) is now
# at row 291, so inserting here pushes it (and everything after) down by one.
$ws.Rows.Item(291).Insert()

$ws.Range("A291").Font.Bold = $true
$ws.Range("A291").HorizontalAlignment = -4108
$ws.Range("A291").VerticalAlignment = -4160
$ws.Range("A291").Borders.LineStyle = 1

$ws.Range("A291").Value = "SAT"
$ws.Range("B291").Value = "San Antonio, United States"
$ws.Range("C291").Value = 29.429461
$ws.Range("D291").Value = -98.487061
$ws.Range("E291").Value = "US"
$ws.Range("F291").Value = "North America"
$ws.Range("G291").Value = "San Antonio"
